$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "24.665.45"
$ws.Range("E2").Value = "  -0.20%  "
Set-TextValue $ws.Range("D3") "1.696.16"
Set-TextValue $ws.Range("D4") "1.003"
$ws.Range("E4").Value = "  +0.20%  "
Set-TextValue $ws.Range("D5") "315.08"
$ws.Range("E5").Value = "  -0.49%  "
Set-TextValue $ws.Range("D6") "1.004"
$ws.Range("E6").Value = "  +0.33%  "
Set-TextValue $ws.Range("D7") "0.3912"
$ws.Range("E7").Value = "  -1.04%  "
Set-TextValue $ws.Range("D8") "0.4056"
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("E10").Value = "  +0.35%  "
Set-TextValue $ws.Range("D11") "53.02"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("E12").Value = "  -1.00%  "
Set-TextValue $ws.Range("D13") "7.659"
$ws.Range("E13").Value = "  +5.64%  "
Set-TextValue $ws.Range("D14") "24.61"
$ws.Range("E14").Value = "  +3.90%  "
Set-TextValue $ws.Range("D15") "0.00001359"
$ws.Range("E15").Value = "  +2.80%  "
Set-TextValue $ws.Range("D16") "7.985"
$ws.Range("E16").Value = "  -0.98%  "
Set-TextValue $ws.Range("D17") "1.689.98"
$ws.Range("E17").Value = "  -0.71%  "
Set-TextValue $ws.Range("D18") "98.56"
Set-TextValue $ws.Range("D19") "0.07112"
$ws.Range("E19").Value = "  +1.18%  "
Set-TextValue $ws.Range("D20") "19.85"
$ws.Range("E20").Value = "  +1.49%  "
Set-TextValue $ws.Range("D21") "7.337"
$ws.Range("E21").Value = "  +4.26%  "
$ws.Range("E22").Value = "  +0.59%  "
Set-TextValue $ws.Range("D23") "14.31"
$ws.Range("E23").Value = "  +0.09%  "
Set-TextValue $ws.Range("D24") "24.631.02"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("E25").Value = "  -7.44%  "
Set-TextValue $ws.Range("D26") "2.352"
$ws.Range("E26").Value = "  -0.46%  "
Set-TextValue $ws.Range("D27") "22.80"
$ws.Range("E27").Value = "  -0.15%  "
Set-TextValue $ws.Range("D28") "162.70"
$ws.Range("E28").Value = "  -0.50%  "
Set-TextValue $ws.Range("D29") "8.467"
$ws.Range("E29").Value = "  +12.25%  "
Set-TextValue $ws.Range("D30") "137.48"
$ws.Range("E30").Value = "  +0.90%  "
Set-TextValue $ws.Range("D31") "5.228"
$ws.Range("E31").Value = "  +0.82%  "
Set-TextValue $ws.Range("D32") "1.877.70"
$ws.Range("E32").Value = "  -0.41%  "
Set-TextValue $ws.Range("D33") "0.08915"
$ws.Range("E33").Value = "  +3.62%  "
Set-TextValue $ws.Range("D34") "7.552"
$ws.Range("E34").Value = "  +5.71%  "
Set-TextValue $ws.Range("D35") "1.049"
$ws.Range("E35").Value = "  -1.95%  "
Set-TextValue $ws.Range("D36") "1.993"
$ws.Range("E36").Value = "  +5.41%  "
$ws.Range("E37").Value = "  +7.95%  "
Set-TextValue $ws.Range("D38") "0.2740"
$ws.Range("E38").Value = "  -0.31%  "
Set-TextValue $ws.Range("D39") "10.78"
$ws.Range("E39").Value = "  -5.45%  "
Set-TextValue $ws.Range("D40") "14.32"
$ws.Range("E40").Value = "  -0.88%  "
Set-TextValue $ws.Range("D41") "0.09126"
$ws.Range("E41").Value = "  -1.24%  "
Set-TextValue $ws.Range("D42") "0.7911"
$ws.Range("E42").Value = "  +3.31%  "
$ws.Range("E43").Value = "  -0.39%  "
Set-TextValue $ws.Range("D44") "16.91"
$ws.Range("E44").Value = "  +5.35%  "
Set-TextValue $ws.Range("D45") "0.7238"
$ws.Range("E45").Value = "  +0.71%  "
$ws.Range("E46").Value = "  +0.06%  "
Set-TextValue $ws.Range("D47") "4.212"
$ws.Range("E47").Value = "  -0.12%  "
Set-TextValue $ws.Range("D48") "1.004"
$ws.Range("E48").Value = "  +0.33%  "
Set-TextValue $ws.Range("D49") "1.333"
$ws.Range("E49").Value = "  +0.39%  "
Set-TextValue $ws.Range("D50") "139.05"
$ws.Range("E50").Value = "  -0.26%  "
Set-TextValue $ws.Range("D51") "91.55"
$ws.Range("E51").Value = "  +1.56%  "
